# No-op: the diff only reorders XML attributes (alphabetical), which
# happens automatically on save; no content change is required.
$d = $word.ActiveDocument
